# Add a new "Obs: 2" textbox to slide 4 (p16), matching the Google Slides
# export that added Shape;298;p16 after the existing shapes.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

$emuPerPt = 12700

# Position / size taken from the target OOXML (EMU -> points).
$x  = 7477125 / $emuPerPt
$y  = 2156550 / $emuPerPt
$cx = 1178700 / $emuPerPt
$cy = 415200  / $emuPerPt

$shp = $s.Shapes.AddTextbox(1, $x, $y, $cx, $cy)
$shp.Name = "Google Shape;298;p16"

# Shape has no fill and no outline.
$shp.Fill.Visible = $false
$shp.Line.Visible = $false

# Text frame insets / wrap / anchor / autosize.
$tf = $shp.TextFrame
$inset = 91425 / $emuPerPt
$tf.MarginLeft = $inset
$tf.MarginRight = $inset
$tf.MarginTop = $inset
$tf.MarginBottom = $inset
$tf.WordWrap = -1
$tf.AutoSize = 0
$tf.VerticalAnchor = 1

$tr = $tf.TextRange
$tr.Text = "Obs: 2"
$tr.LanguageID = "pt-BR"

# Paragraph formatting: left aligned, no space before/after, no bullet.
$pf = $tr.ParagraphFormat
$pf.Alignment = 1
$pf.SpaceBefore = 0
$pf.SpaceAfter = 0
$pf.Bullet.Visible = 0

# Font formatting.
$f = $tr.Font
$f.Name = "Nunito"
$f.NameFarEast = "Nunito"
$f.NameOther = "Nunito"
$f.Size = 13
$f.Color.ObjectThemeColor = 3
